# Update: Pra LOP detail + Form P1 + Print Form P1
#
# Three tables in the document are touched:
#   Table 2 (COM index) - "No/Lokasi/Item/Qty/Unit" small table: two columns get
#     resized by a single dxa (639->638, 1401->1402).
#   Table 3 (COM index) - the "Pra LOP" detail table with ${harga_mrc}/${harga_otc}
#     placeholders: header text is relabeled (Total Harga (MRC) -> OTC,
#     Total Harga OTC -> Total Harga), the placeholder values shift one column to
#     the left, and the last two columns are resized (1511->1513, 1929->1927).
#   Table 4 (COM index) - a twin of table 3 used elsewhere in the form: only the
#     same two columns are resized (1511->1513, 1929->1927); no text changes.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Table 2 : small "No / Lokasi / Item / Qty / Unit" table - column widths only
# ---------------------------------------------------------------------------
$t2 = $d.Tables.Item(2)
$t2.Rows.Item(1).Cells.Item(4).Width = 31.9    # 639 dxa -> 638 dxa
$t2.Rows.Item(1).Cells.Item(5).Width = 70.1    # 1401 dxa -> 1402 dxa

# ---------------------------------------------------------------------------
# Table 3 : Pra LOP detail table with the harga_mrc / harga_otc placeholders
# ---------------------------------------------------------------------------
$t3 = $d.Tables.Item(3)

# Header row: relabel the last two header cells.
$t3.Rows.Item(1).Cells.Item(8).Range.Text = "OTC"
$t3.Rows.Item(1).Cells.Item(9).Range.Text = "Total Harga"

# Data row: the two placeholder values move one cell to the left, the former
# last cell becomes empty.
$t3.Rows.Item(2).Cells.Item(7).Range.Text = "`${harga_mrc}"
$t3.Rows.Item(2).Cells.Item(8).Range.Text = "`${harga_otc}"
$t3.Rows.Item(2).Cells.Item(9).Range.Text = ""

# Column widths for the header + data rows (1511 -> 1513, 1929 -> 1927).
$t3.Rows.Item(1).Cells.Item(8).Width = 75.65
$t3.Rows.Item(1).Cells.Item(9).Width = 96.35
$t3.Rows.Item(2).Cells.Item(8).Width = 75.65
$t3.Rows.Item(2).Cells.Item(9).Width = 96.35

# ---------------------------------------------------------------------------
# Table 4 : twin table elsewhere in the form - column widths only
# ---------------------------------------------------------------------------
$t4 = $d.Tables.Item(4)
$t4.Rows.Item(1).Cells.Item(8).Width = 75.65
$t4.Rows.Item(1).Cells.Item(9).Width = 96.35
$t4.Rows.Item(2).Cells.Item(8).Width = 75.65
$t4.Rows.Item(2).Cells.Item(9).Width = 96.35
